$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains price text that can look like plain numbers (e.g. "1.00", "0.999").
# Force the whole Price column to Text format before writing so Excel keeps the values
# as strings (preserving formatting such as trailing zeros) instead of silently
# converting them to numeric cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "97.590.74"
$ws.Range("E2").Value = "  +1.08%  "
$ws.Range("D3").Value = "3.728.22"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("E5").Value = "  +13.13%  "
$ws.Range("D6").Value = "238.48"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "657.67"
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  +5.18%  "
$ws.Range("E9").Value = "  +3.85%  "
$ws.Range("D10").Value = "0.999"
$ws.Range("E10").Value = "  -0.05%  "
$ws.Range("D11").Value = "3.724.21"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("E12").Value = "  +16.97%  "
$ws.Range("D13").Value = "44.84"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").Value = "4.423.83"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "97.377.10"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "9.26"
$ws.Range("E18").Value = "  +2.42%  "
$ws.Range("D19").Value = "3.716.34"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "13.11"
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").Value = "18.89"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").Value = "529.04"
$ws.Range("E23").Value = "  +1.27%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  +10.10%  "
$ws.Range("D26").Value = "117.51"
$ws.Range("E26").Value = "  +14.99%  "
$ws.Range("D27").Value = "6.91"
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").Value = "0.212"
$ws.Range("E28").Value = "  +25.91%  "
$ws.Range("D29").Value = "13.43"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "0.600"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "641.40"
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +4.28%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "6.83"
$ws.Range("E42").Value = "  -3.81%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.496"
$ws.Range("E43").Value = "  +10.74%  "
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("D46").Value = "0.970"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("D49").Value = "8.81"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  +4.57%  "

# Restore the default cell style on the Price column so the on-disk formatting
# matches the original workbook (values remain text).
$ws.Range("D2:D51").Style = "Normal"

